{"js": "// \"Remove indentation from compact style\"\n//\n// The custom paragraph style \"Compact\" (based on \"Body Text\") inherits a\n// first-line indent from its base style. This removes that inherited\n// indentation by setting the style's own first-line indent to 0 (which\n// OOXML serializes as <w:ind w:firstLine=\"0\"/> inside the style's <w:pPr>).\n\nconst styles = context.document.getStyles();\nconst compact = styles.getByName(\"Compact\");\n\n// Zeroing the first-line indent directly on the style definition.\ncompact.paragraphFormat.firstLineIndent = 0;\n\nawait context.sync();\n", "ps1": "# \"Remove indentation from compact style\"\n#\n# The custom paragraph style \"Compact\" (based on \"Body Text\") inherits a\n# first-line indent from its base style. This removes that inherited\n# indentation by setting the style's own first-line indent to 0 (which\n# OOXML serializes as <w:ind w:firstLine=\"0\"/> inside the style's <w:pPr>).\n\n$d = $word.ActiveDocument\n$compact = $d.Styles(\"Compact\")\n\n# Zeroing the first-line indent directly on the style definition.\n$compact.ParagraphFormat.FirstLineIndent = 0\n"}
